$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "56.836.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "  +2.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "2.491.89"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = "  -1.13%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value2 = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "491.64"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "152.23"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = "  +7.40%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = "  +0.12%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = "  -0.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "2.503.32"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value2 = "  -0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "5.73"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = "  +3.08%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = "  -0.88%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = "  +0.14%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = "  +0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "2.925.35"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = "  -0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "56.729.91"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "  +1.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "21.21"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value2 = "  +1.14%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = "  -2.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "2.497.39"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = "  -0.35%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value2 = "  +3.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "10.32"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value2 = "  +2.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "321.49"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = "  -0.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.997"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = "  -0.15%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = "  +1.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "58.43"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = "  +0.09%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value2 = "  -0.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = "  +0.32%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value2 = "  -6.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "2.595.29"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value2 = "  -0.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "7.57"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value2 = "  +0.96%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value2 = "  -0.25%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value2 = "  +0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "151.82"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value2 = "  +0.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "18.33"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value2 = "  -0.22%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value2 = "  +0.44%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value2 = "  +0.07%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value2 = "  +1.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "3.77"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value2 = "  +0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.868"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value2 = "  -1.30%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value2 = "  +4.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "34.19"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value2 = "  -1.06%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value2 = "  +1.98%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value2 = "  +0.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.615"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value2 = "  -0.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.997"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value2 = "  +0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "267.19"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = "  +2.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "4.79"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value2 = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.0933"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value2 = "  +1.74%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value2 = "VeChain"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.0229"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value2 = "  +0.78%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value2 = "WhiteBITCoin"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "10.21"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value2 = "  +0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "17.85"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = "  +1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "1.895.27"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value2 = "  -5.89%  "
